# Add Inventory (Bulk) sheet format update
# - Reworks the header row (row 1) columns/labels
# - Adds new columns (e-File Number, Price, Section Of Center)
# - Restyles the placeholder data row (row 2)
# - Adjusts column widths to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) values
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Item Category*`n(Make sure to`nwrite only the`ncategories`npresent in`ndatabase)"
$ws.Range("B1").Value = "Manufacturer*`n(Please enter`nonly`nmanufacturers`nin the `ndatabase)"
$ws.Range("C1").Value = "Item Type`n(New/Old)`nDefault: New"
$ws.Range("D1").Value = "Date of Purchase*"
$ws.Range("E1").Value = "Bill in which`n charged"
$ws.Range("F1").Value = "Model*"
$ws.Range("G1").Value = "Serial No.*"
$ws.Range("H1").Value = "Specifications`n(ComputerNumber)"
$ws.Range("I1").Value = "Warranty`n(Only write one of`nthe following:`nWarranty`nAMC`nNo Warranty / No AMC)`nNo Warranty / No AMC`nby default"
$ws.Range("J1").Value = "Warranty Expiry"
$ws.Range("K1").Value = "Vendor"
$ws.Range("L1").Value = "e-File Number"
$ws.Range("M1").Value = "Price"
$ws.Range("N1").Value = "Section Of Center"

# ---------------------------------------------------------------------------
# 2. Header row formatting: bold, Arial, automatic/theme text color, top
#    aligned (matches the existing header look for the whole A1:N1 range).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:N1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 1
$headerRange.VerticalAlignment = -4160
$headerRange.NumberFormat = "General"

# Item Type column keeps a text number format (as Quantity* used to)
$ws.Range("C1").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3. Placeholder data row (row 2) formatting: regular weight, Arial, themed
#    text color. Date of Purchase* now lives in column D, so the date number
#    format follows it there.
# ---------------------------------------------------------------------------
$row2Range = $ws.Range("A2:C2,E2:K2,M2")
$row2Range.Font.Name = "Arial"
$row2Range.Font.Bold = $false
$row2Range.Font.ThemeColor = 1
$row2Range.NumberFormat = "General"

$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Bold = $false
$ws.Range("D2").Font.ThemeColor = 1
$ws.Range("D2").NumberFormat = "m/d/yyyy"

# Column L no longer carries the old "Vendor" placeholder formatting
$ws.Range("L2").ClearFormats()

# ---------------------------------------------------------------------------
# 4. Column widths for the new layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14.43  # B back to default width
$ws.Columns.Item(3).ColumnWidth = 16.86  # C
$ws.Columns.Item(4).ColumnWidth = 21.0   # D
$ws.Columns.Item(5).ColumnWidth = 17.71  # E
$ws.Columns.Item(6).ColumnWidth = 27.43  # F
$ws.Columns.Item(7).ColumnWidth = 21.43  # G
$ws.Columns.Item(8).ColumnWidth = 24.43  # H
$ws.Columns.Item(9).ColumnWidth = 22.43  # I
$ws.Columns.Item(10).ColumnWidth = 16.29 # J
$ws.Columns.Item(11).ColumnWidth = 14.43 # K back to default width
$ws.Columns.Item(14).ColumnWidth = 18.43 # N
